$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("compra")

# Fix the broken linkage between "compra" purchase-order rows and the
# "requi" requisition rows: rows 21-30 (items 6-15) were either blank or
# pointing at the wrong requi rows; row 31 (item 16) incorrectly carried
# formulas that belong nowhere. Re-point rows 21-30 sequentially at
# requi rows 25-34, and clear out row 31.

for ($i = 0; $i -le 9; $i++) {
    $compraRow = 21 + $i
    $requiRow = 25 + $i

    $ws.Range("C$compraRow").Formula = "=requi!J$requiRow"
    $ws.Range("D$compraRow").Formula = "=requi!I$requiRow"
    $ws.Range("E$compraRow").Formula = "=requi!C$requiRow"
    $ws.Range("F$compraRow").Formula = "=requi!D$requiRow"
}

# Row 31 (item 16) no longer references requi!36 - clear it out.
$ws.Range("C31").Formula = ""
$ws.Range("D31").Formula = ""
$ws.Range("E31").Formula = ""
$ws.Range("F31").Formula = ""

# Restore the selection left on the "compra" sheet to F20:K20.
[void]$ws.Range("F20:K20").Select()
